# Auto-generated edit script applying numeric corrections per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1549
$ws.Range("J17").Value = 1549
$ws.Range("L17").Value = 4647
$ws.Range("N17").Value = -4983
$ws.Range("H62").Value = 4099.5
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 4099.5
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H101").Value = 3374.875
$ws.Range("I101").Value = 601.6
$ws.Range("K101").Value = 1804.8
$ws.Range("M101").Value = -182.8000000000002
$ws.Range("H106").Value = 83337580
$ws.Range("I106").Value = 166667660
$ws.Range("K106").Value = 166667660
$ws.Range("M106").Value = -166667029
$ws.Range("H116").Value = 41055892
$ws.Range("I116").Value = 86915890
$ws.Range("K116").Value = 86915890
$ws.Range("M116").Value = -86912448
$ws.Range("H138").Value = 7165.811
$ws.Range("I138").Value = 1617.7142
$ws.Range("J138").Value = 8460.366
$ws.Range("K138").Value = 4853.142599999999
$ws.Range("L138").Value = 25381.098
$ws.Range("M138").Value = 286.8574000000008
$ws.Range("N138").Value = -35661.098

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1343847.9
$ws.Range("I2").Value = 1940062.5
$ws.Range("K2").Value = 1940062.5
$ws.Range("M2").Value = -1939949.5
$ws.Range("H45").Value = 4299.6665
$ws.Range("I45").Value = 4299.6665
$ws.Range("K45").Value = 4299.6665
$ws.Range("M45").Value = -3922.6665
$ws.Range("H61").Value = 13282.895
$ws.Range("J61").Value = 15971.8
$ws.Range("L61").Value = 15971.8
$ws.Range("N61").Value = -16395.8
$ws.Range("H110").Value = 6883.75
$ws.Range("I110").Value = 2845
$ws.Range("J110").Value = 19000
$ws.Range("K110").Value = 2845
$ws.Range("L110").Value = 19000
$ws.Range("M110").Value = -800
$ws.Range("N110").Value = -23090
$ws.Range("H116").Value = 1343847.9
$ws.Range("I116").Value = 1940062.5
$ws.Range("K116").Value = 1940062.5
$ws.Range("M116").Value = -1937768.5
$ws.Range("H122").Value = 8996.583000000001
$ws.Range("I122").Value = 7708.7144
$ws.Range("K122").Value = 23126.1432
$ws.Range("M122").Value = -20676.1432
$ws.Range("H132").Value = 15509.935
$ws.Range("I132").Value = 15225.878
$ws.Range("K132").Value = 45677.63400000001
$ws.Range("M132").Value = -43147.63400000001
$ws.Range("H136").Value = 13282.895
$ws.Range("J136").Value = 15971.8
$ws.Range("L136").Value = 47915.39999999999
$ws.Range("N136").Value = -53015.39999999999
$ws.Range("H140").Value = 113133.336
$ws.Range("J140").Value = 113133.336
$ws.Range("L140").Value = 113133.336
$ws.Range("N140").Value = -123493.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1343847.9
$ws.Range("I3").Value = 1940062.5
$ws.Range("K3").Value = 1940062.5
$ws.Range("M3").Value = -1939948.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1433.9333
$ws.Range("I16").Value = 1330.8182
$ws.Range("J16").Value = 1717.5
$ws.Range("K16").Value = 1330.8182
$ws.Range("L16").Value = 1717.5
$ws.Range("M16").Value = -1043.8182
$ws.Range("N16").Value = -2291.5
$ws.Range("H31").Value = 2069.2632
$ws.Range("I31").Value = 1055.2333
$ws.Range("J31").Value = 5871.875
$ws.Range("K31").Value = 1055.2333
$ws.Range("L31").Value = 5871.875
$ws.Range("M31").Value = -760.2333000000001
$ws.Range("N31").Value = -6461.875
$ws.Range("H34").Value = 2069.2632
$ws.Range("I34").Value = 1055.2333
$ws.Range("J34").Value = 5871.875
$ws.Range("K34").Value = 1055.2333
$ws.Range("L34").Value = 5871.875
$ws.Range("M34").Value = -853.2333000000001
$ws.Range("N34").Value = -6275.875
$ws.Range("H58").Value = 401818.8
$ws.Range("I58").Value = 528048.6
$ws.Range("K58").Value = 528048.6
$ws.Range("M58").Value = -527845.6
$ws.Range("H99").Value = 7150.1875
$ws.Range("I99").Value = 5501.6
$ws.Range("J99").Value = 7899.5454
$ws.Range("K99").Value = 5501.6
$ws.Range("L99").Value = 7899.5454
$ws.Range("M99").Value = -4003.6
$ws.Range("N99").Value = -10895.5454
$ws.Range("H113").Value = 1433.9333
$ws.Range("I113").Value = 1330.8182
$ws.Range("J113").Value = 1717.5
$ws.Range("K113").Value = 1330.8182
$ws.Range("L113").Value = 1717.5
$ws.Range("M113").Value = 839.1818000000001
$ws.Range("N113").Value = -6057.5
$ws.Range("H126").Value = 7150.1875
$ws.Range("I126").Value = 5501.6
$ws.Range("J126").Value = 7899.5454
$ws.Range("K126").Value = 16504.8
$ws.Range("L126").Value = 23698.6362
$ws.Range("M126").Value = -14034.8
$ws.Range("N126").Value = -28638.6362
$ws.Range("H134").Value = 2333.932
$ws.Range("I134").Value = 2304.7144
$ws.Range("J134").Value = 2947.5
$ws.Range("K134").Value = 6914.1432
$ws.Range("L134").Value = 8842.5
$ws.Range("M134").Value = -4379.1432
$ws.Range("N134").Value = -13912.5
$ws.Range("H136").Value = 401818.8
$ws.Range("I136").Value = 528048.6
$ws.Range("K136").Value = 1584145.8
$ws.Range("M136").Value = -1581595.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 381.9091
$ws.Range("J12").Value = 393.2857
$ws.Range("L12").Value = 1179.8571
$ws.Range("N12").Value = -1525.8571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 650088.8
$ws.Range("I122").Value = 690581.9
$ws.Range("K122").Value = 2071745.7
$ws.Range("M122").Value = -2069295.7
$ws.Range("H126").Value = 4453.4287
$ws.Range("I126").Value = 2969.5
$ws.Range("K126").Value = 8908.5
$ws.Range("M126").Value = -6438.5
$ws.Range("H132").Value = 4088.5715
$ws.Range("I132").Value = 4101.645
$ws.Range("K132").Value = 12304.935
$ws.Range("M132").Value = -9774.935000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7273.4287
$ws.Range("I40").Value = 7000.75
$ws.Range("J40").Value = 7637
$ws.Range("K40").Value = 7000.75
$ws.Range("L40").Value = 7637
$ws.Range("M40").Value = -6864.75
$ws.Range("N40").Value = -7909
$ws.Range("H127").Value = 250104900
$ws.Range("J127").Value = 139865.33
$ws.Range("L127").Value = 139865.33
$ws.Range("N127").Value = -149785.33
$ws.Range("H132").Value = 4097.508
$ws.Range("I132").Value = 3165.848
$ws.Range("J132").Value = 6353.1055
$ws.Range("K132").Value = 9497.544
$ws.Range("L132").Value = 19059.3165
$ws.Range("M132").Value = -6967.544
$ws.Range("N132").Value = -24119.3165

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 44999
$ws.Range("J43").Value = 44999
$ws.Range("L43").Value = 44999
$ws.Range("N43").Value = -45297
$ws.Range("H122").Value = 4313.5557
$ws.Range("I122").Value = 4098.7085
$ws.Range("J122").Value = 6032.3335
$ws.Range("K122").Value = 12296.1255
$ws.Range("L122").Value = 18097.0005
$ws.Range("M122").Value = -9846.125499999998
$ws.Range("N122").Value = -22997.0005
$ws.Range("H126").Value = 2750
$ws.Range("I126").Value = 2750
$ws.Range("K126").Value = 8250
$ws.Range("M126").Value = -5780
$ws.Range("H132").Value = 12628519
$ws.Range("I132").Value = 1793870.9
$ws.Range("J132").Value = 38464988
$ws.Range("K132").Value = 5381612.699999999
$ws.Range("L132").Value = 115394964
$ws.Range("M132").Value = -5379082.699999999
$ws.Range("N132").Value = -115400024
